$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.048.39"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "3.207.54"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'602.41"
$ws.Range("E5").Value = "  +4.06%  "
$ws.Range("D6").Value = "'153.88"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.208.43"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'6.15"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "'38.59"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "3.732.21"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "66.198.45"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("E17").Value = "  +3.83%  "
$ws.Range("D18").Value = "3.208.78"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "'511.07"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("E21").Value = "  +4.85%  "
$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").Value = "'15.20"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'85.14"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'3.01"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("D31").Value = "'6.88"
$ws.Range("E31").Value = "  +8.77%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "'6.62"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'55.31"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'0.0911"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "'484.00"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "'0.0420"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "'2.97"
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "'0.298"
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("D43").Value = "'0.119"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("D45").Value = "2.961.36"
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "0.0₃0641"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").Value = "'28.83"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("D51").Value = "'34.08"
$ws.Range("E51").Value = "  +5.53%  "
